# repos_data.xlsx tracked a list of repo URLs in column A of "Feuil1".
# Row 2 contained https://github.com/gautamvr/DocumentProcessor_GCP.git
# (a project that was not using the batch API) and row 3 contained
# https://github.com/ovokpus/Python-Azure-AI-REST-APIs.git.
#
# This edit ("fix not using batch api") removes the gautamvr entry and
# keeps only the ovokpus repo, also cleaning up its URL by dropping the
# trailing ".git" suffix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrl = "https://github.com/ovokpus/Python-Azure-AI-REST-APIs"

# Re-point the existing A2 hyperlink (currently the gautamvr/.git repo) at
# the cleaned-up ovokpus URL, reusing the cell/hyperlink/style already in
# place rather than creating a brand new one.
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Address -eq "https://github.com/gautamvr/DocumentProcessor_GCP.git") {
        $hl.Address = $newUrl
    }
}
$ws.Range("A2").Value = $newUrl

# Remove the hyperlink that used to sit on A3 (the old, ".git"-suffixed
# ovokpus link) and then delete that now-redundant row entirely.
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Address -eq "https://github.com/ovokpus/Python-Azure-AI-REST-APIs.git") {
        $hl.Delete()
    }
}
$ws.Range("A3").EntireRow.Delete()

# Match the saved selection state (cursor parked on the now-empty A3).
$ws.Range("A3").Select()
